# The deck originally used the "Integral" (Red Violet) design for its
# slides/slide master (ppt/theme/theme1.xml) and carried a separate,
# unused "Office Theme" colour set for the notes master
# (ppt/theme/theme2.xml). The edit swaps the two theme colour schemes so
# that the slide master now uses the default Office palette.
#
# PowerPoint exposes the 12 theme colour slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink - in that order) through
# Slide.ThemeColorScheme.Colors(1..12).RGB, which writes straight back
# into the presentation's master theme part.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette = the standard Office theme colours.
$officeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
